$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Helper: replace a paragraph's run content (everything except the trailing
# paragraph mark, so the paragraph's <w:pPr> is left untouched) with a raw
# WordprocessingML fragment.
# ---------------------------------------------------------------------------
function Set-ParaInnerXml($para, [string]$innerXml) {
    $start = $para.Range.Start
    $end = $para.Range.End - 1
    $r = $d.Range($start, $end)
    $r.InsertXML("<w:p $wns>" + $innerXml + "</w:p>") | Out-Null
}

# ---------------------------------------------------------------------------
# 1. "Oleh :" paragraph (index 5) - mark "Oleh" as a spelling item and wrap
#    the whole "Oleh :" fragment in a grammar-check range.
# ---------------------------------------------------------------------------
$pOleh = $d.Paragraphs.Item(5)
$xmlOleh = '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Oleh</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>'
Set-ParaInnerXml $pOleh $xmlOleh

# ---------------------------------------------------------------------------
# 2. "Ayu Permata Sari" paragraph (index 6) - split into Ayu / Permata Sari
#    words, each flagged by the spell checker.
# ---------------------------------------------------------------------------
$pAyu = $d.Paragraphs.Item(6)
$rprAyu = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$xmlAyu = '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprAyu<w:t>Ayu</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$rprAyu<w:t xml:space=`"preserve`"> </w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprAyu<w:t>Permata</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$rprAyu<w:t xml:space=`"preserve`"> Sari</w:t></w:r>"
Set-ParaInnerXml $pAyu $xmlAyu

# ---------------------------------------------------------------------------
# 3. "Muhammad Nur Ikhsan" paragraph (index 9) - split into Muhammad / Nur /
#    Ikhsan, with "Nur" and "Ikhsan" flagged by the spell checker.
# ---------------------------------------------------------------------------
$pMuh = $d.Paragraphs.Item(9)
$rprMuh = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$xmlMuh = "<w:r>$rprMuh<w:t xml:space=`"preserve`">Muhammad </w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprMuh<w:t>Nur</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$rprMuh<w:t xml:space=`"preserve`"> </w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprMuh<w:t>Ikhsan</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>'
Set-ParaInnerXml $pMuh $xmlMuh

# ---------------------------------------------------------------------------
# 4. "Nur Rahmawati" paragraph (index 10) - split into Nur / Rahmawati, both
#    flagged by the spell checker, then add a new blank paragraph after it.
# ---------------------------------------------------------------------------
$pNur = $d.Paragraphs.Item(10)
$rprNur = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$xmlNur = '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprNur<w:t>Nur</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$rprNur<w:t xml:space=`"preserve`"> </w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprNur<w:t>Rahmawati</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>'
Set-ParaInnerXml $pNur $xmlNur

# Insert a new, totally empty paragraph right after "Nur Rahmawati".
$pNur = $d.Paragraphs.Item(10)
$insertionPoint = $d.Range($pNur.Range.End, $pNur.Range.End)
$insertionPoint.InsertXML("<w:p $wns/>") | Out-Null

# ---------------------------------------------------------------------------
# 5. "Andi syahjaratu " paragraph (now index 12) - capitalise to
#    "Andi Syahjaratu ", split into Andi / Syahjaratu, both flagged by the
#    spell checker.
# ---------------------------------------------------------------------------
$pAndi = $d.Paragraphs.Item(12)
$rprAndi = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$xmlAndi = '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprAndi<w:t>Andi</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$rprAndi<w:t xml:space=`"preserve`"> </w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rprAndi<w:t>Syahjaratu</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$rprAndi<w:t xml:space=`"preserve`"> </w:t></w:r>"
Set-ParaInnerXml $pAndi $xmlAndi
